# Read Test data from Excel
# Populate Sheet1 with a "URL" header and a hyperlinked Google URL, matching
# what a user would get after typing the values in and letting Excel
# auto-format the link (Hyperlink cell style: underlined, theme color 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: plain text header
$ws.Range("A1").Value = "URL"

# B1: the URL text itself
$ws.Range("B1").Value = "https://www.google.com"

# Turn B1 into a real hyperlink pointing at the URL (this also applies the
# built-in "Hyperlink" cell style: underline + theme color 10).
$link = $ws.Hyperlinks.Add($ws.Range("B1"), "")
$link.Address = "https://www.google.com"

# Leave the selection where it would land after entering the URL and
# pressing Enter, i.e. one row below the link.
$ws.Range("B2").Select() | Out-Null
